# draft plot C.i for mitigation, and make notes for further investigation
# (absence of salt marsh, MRE x coral reef in mitigation outcomes)
#
# Adds two new lookup blocks to the factor_aesthetics sheet:
#   - ecosystem_type (rows 23-26): Salt marsh / Mangrove / Seagrass / Coral reef
#   - adapt_to_threat (rows 27-28): Human / Natural
#
# Cell values are written in the exact order needed so the shared-string
# table is built up in the same sequence as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ecosystem_type block (rows 23-26) ------------------------------------

# variable column
$ws.Range("A23").Value = "ecosystem_type"
$ws.Range("A24").Value = "ecosystem_type"
$ws.Range("A25").Value = "ecosystem_type"
$ws.Range("A26").Value = "ecosystem_type"

# level column (note: row 26 filled before row 25, matching source order)
$ws.Range("B23").Value = "Salt_marsh"
$ws.Range("B24").Value = "Mangrove"
$ws.Range("B26").Value = "Coral_reef"
$ws.Range("B25").Value = "Seagrass"

# label column
$ws.Range("C23").Value = "Salt marsh"
$ws.Range("C24").Value = "Mangrove"
$ws.Range("C26").Value = "Coral reef"
$ws.Range("C25").Value = "Seagrass"

# order column
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("D26").Value = 4

# colour column
$ws.Range("E23").Value = "#E69F00"
$ws.Range("E24").Value = "#009E73"
$ws.Range("E25").Value = "#56B4E9"
$ws.Range("E26").Value = "#D55E00"

# --- adapt_to_threat block (rows 27-28) -----------------------------------

$ws.Range("A27").Value = "adapt_to_threat"
$ws.Range("A28").Value = "adapt_to_threat"

$ws.Range("B27").Value = "Human"
$ws.Range("B28").Value = "Natural"

$ws.Range("C27").Value = "Human"
$ws.Range("C28").Value = "Nature"

$ws.Range("D27").Value = 1
$ws.Range("D28").Value = 2

# --- view state: scroll down and select the next empty row ----------------

$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select()
